$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.908.73"
$ws.Range("E2").Value = "  +2.58%  "
$ws.Range("D3").Value = "1.905.72"
$ws.Range("E3").Value = "  +1.01%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'246.78"
$ws.Range("E5").Value = "  +1.33%  "
$ws.Range("D6").Value = "'1.0000"
$ws.Range("E6").Value = "  +0.09%  "
$ws.Range("E7").Value = "  +0.62%  "
$ws.Range("D8").Value = "'0.3004"
$ws.Range("E8").Value = "  +2.21%  "
$ws.Range("D9").Value = "'0.06867"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").Value = "1.906.69"
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  +3.64%  "
$ws.Range("D12").Value = "'0.07352"
$ws.Range("E12").Value = "  +2.42%  "
$ws.Range("D13").Value = "'91.92"
$ws.Range("E13").Value = "  +6.73%  "
$ws.Range("D14").Value = "'5.139"
$ws.Range("E14").Value = "  +5.80%  "
$ws.Range("D15").Value = "'0.6846"
$ws.Range("E15").Value = "  +2.87%  "
$ws.Range("D16").Value = "30.889.88"
$ws.Range("E16").Value = "  +2.63%  "
$ws.Range("D17").Value = "'0.000008096"
$ws.Range("E17").Value = "  +1.86%  "
$ws.Range("D18").Value = "'13.47"
$ws.Range("E18").Value = "  +5.36%  "
$ws.Range("D19").Value = "'0.9995"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "2.154.11"
$ws.Range("E20").Value = "  +1.53%  "
$ws.Range("E21").Value = "  +0.22%  "
$ws.Range("D22").Value = "'4.888"
$ws.Range("E22").Value = "  +2.60%  "
$ws.Range("D23").Value = "'183.20"
$ws.Range("E23").Value = "  +35.94%  "
$ws.Range("D24").Value = "'6.132"
$ws.Range("E24").Value = "  +9.36%  "
$ws.Range("D25").Value = "'9.403"
$ws.Range("E25").Value = "  +2.82%  "
$ws.Range("D26").Value = "'154.23"
$ws.Range("E26").Value = "  +1.85%  "
$ws.Range("D27").Value = "'18.75"
$ws.Range("E27").Value = "  +11.53%  "
$ws.Range("D28").Value = "'1.961"
$ws.Range("E28").Value = "  +2.19%  "
$ws.Range("E29").Value = "  +1.32%  "
$ws.Range("D30").Value = "'4.407"
$ws.Range("E30").Value = "  +5.74%  "
$ws.Range("D31").Value = "'0.08992"
$ws.Range("E31").Value = "  +3.42%  "
$ws.Range("D32").Value = "'4.085"
$ws.Range("E32").Value = "  +3.40%  "
$ws.Range("D33").Value = "'0.05312"
$ws.Range("E33").Value = "  +6.06%  "
$ws.Range("D34").Value = "'0.7503"
$ws.Range("E34").Value = "  +6.11%  "
$ws.Range("D35").Value = "'1.147"
$ws.Range("D36").Value = "'2.699"
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").Value = "'0.01924"
$ws.Range("E37").Value = "  +16.41%  "
$ws.Range("D38").Value = "'2.735"
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("E39").Value = "  +0.23%  "
$ws.Range("D40").Value = "'0.9410"
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("D41").Value = "'0.4410"
$ws.Range("E41").Value = "  +4.85%  "
$ws.Range("D42").Value = "'106.47"
$ws.Range("E42").Value = "  +3.87%  "
$ws.Range("D43").Value = "'5.865"
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  +0.06%  "
$ws.Range("E45").Value = "  +4.26%  "
$ws.Range("D46").Value = "'0.1363"
$ws.Range("E46").Value = "  +8.25%  "
$ws.Range("E47").Value = "  +2.82%  "
$ws.Range("E48").Value = "  +5.65%  "
$ws.Range("D49").Value = "'8.629"
$ws.Range("E49").Value = "  +4.92%  "
$ws.Range("D50").Value = "'33.63"
$ws.Range("E50").Value = "  +3.57%  "
$ws.Range("E51").Value = "  +3.96%  "
